$wb = $excel.ActiveWorkbook

# --- Sheet 1 "Calibração": B2 becomes the text "s" instead of the number 1.38
$calibracao = $wb.Worksheets.Item("Calibração")
$calibracao.Range("B2").Value = "s"
# Update the selection on this (no longer active) sheet
$null = $calibracao.Range("A2:A26").Select()

# --- Add a new sheet "Calibração com " at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$novaCalibracao = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$novaCalibracao.Name = "Calibração com "

$novaCalibracao.Range("A1:D1").Font.Bold = $true
$novaCalibracao.Range("A1").Value = "Ângulo (graus)"
$novaCalibracao.Range("B1").Value = "Tensão (V)"
$novaCalibracao.Range("C1").Value = "Tensão (V)"
$null = $novaCalibracao.Range("A2").Select()

# --- Sheet 3 "Formas de onda" becomes the active sheet/tab, with a new selection
$formas = $wb.Worksheets.Item("Formas de onda")
$formas.Activate()
$null = $formas.Range("C7").Select()
